# Apply "added control mode switch" edits to the BOM worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities affected by adding a control-mode switch.
$ws.Range("D11").Value = 5
$ws.Range("D12").Value = 5
$ws.Range("D18").Value = 3

# Update the current view/selection state on the sheet.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D13").Select()
